$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task rows to append below the existing data (rows 1-3 already present
# with Description / Due Date / Due Time columns).
$rows = @(
    @("drop the kids in school", "2024-11-18", "11:00"),
    @("kids holiday", "2024-06-28", "11:40"),
    @("dinner with friends", "2024-08-13", "08:00"),
    @("school holidays", "2024-07-12", "11:30")
)

$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]

    # Due Date / Due Time look like dates/times; force them to be stored as
    # plain text (matching the rest of the column) instead of letting Excel
    # auto-convert them to date/time serial numbers. Setting the format to
    # Text ("@") before assigning keeps the literal string, and clearing the
    # format afterwards avoids leaving a custom number-format style behind
    # on the cell (matches the General-styled cells elsewhere in the sheet).
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 2).ClearFormats()

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 3).ClearFormats()
}
